$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the test-result flags (PASS/FAIL) in the Results column for rows 2-4
$ws.Range("E2:E4").ClearContents()

# Update the selection to match the post-edit state
$ws.Range("E2:E4").Select()
